# Updates market-price-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across all eight crafting-class leve sheets, refreshing the scraped market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2124.0454
$ws.Range("J17").Value = 2158.524
$ws.Range("L17").Value = 6475.572
$ws.Range("N17").Value = -6811.572
$ws.Range("H26").Value = 16500
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20688
$ws.Range("H74").Value = 88935.05499999999
$ws.Range("I74").Value = 97692.17999999999
$ws.Range("J74").Value = 14499.5
$ws.Range("K74").Value = 97692.17999999999
$ws.Range("L74").Value = 14499.5
$ws.Range("M74").Value = -96756.17999999999
$ws.Range("N74").Value = -16371.5
$ws.Range("H75").Value = 39998.332
$ws.Range("J75").Value = 39998.332
$ws.Range("L75").Value = 39998.332
$ws.Range("N75").Value = -41870.332
$ws.Range("H77").Value = 88935.05499999999
$ws.Range("I77").Value = 97692.17999999999
$ws.Range("J77").Value = 14499.5
$ws.Range("K77").Value = 488460.9
$ws.Range("L77").Value = 72497.5
$ws.Range("M77").Value = -483780.9
$ws.Range("N77").Value = -81857.5
$ws.Range("H78").Value = 39998.332
$ws.Range("J78").Value = 39998.332
$ws.Range("L78").Value = 119994.996
$ws.Range("N78").Value = -129354.996
$ws.Range("H80").Value = 5653.273
$ws.Range("I80").Value = 2705.3333
$ws.Range("J80").Value = 9190.799999999999
$ws.Range("K80").Value = 8115.999899999999
$ws.Range("L80").Value = 27572.4
$ws.Range("M80").Value = -7117.999899999999
$ws.Range("N80").Value = -29568.4
$ws.Range("H83").Value = 5653.273
$ws.Range("I83").Value = 2705.3333
$ws.Range("J83").Value = 9190.799999999999
$ws.Range("K83").Value = 24347.9997
$ws.Range("L83").Value = 82717.2
$ws.Range("M83").Value = -19355.9997
$ws.Range("N83").Value = -92701.2
$ws.Range("H92").Value = 566.8889
$ws.Range("J92").Value = 765.6667
$ws.Range("L92").Value = 765.6667
$ws.Range("N92").Value = -3261.6667
$ws.Range("H101").Value = 25000444
$ws.Range("I101").Value = 25000444
$ws.Range("K101").Value = 75001332
$ws.Range("M101").Value = -74999710
$ws.Range("H112").Value = 3318.75
$ws.Range("J112").Value = 3488
$ws.Range("L112").Value = 10464
$ws.Range("N112").Value = -12680
$ws.Range("H116").Value = 5088
$ws.Range("I116").Value = 2598.6
$ws.Range("K116").Value = 2598.6
$ws.Range("M116").Value = 843.4000000000001
$ws.Range("H129").Value = 3832.8
$ws.Range("I129").Value = 4031.7144
$ws.Range("J129").Value = 3368.6667
$ws.Range("K129").Value = 12095.1432
$ws.Range("L129").Value = 10106.0001
$ws.Range("M129").Value = -7095.143199999999
$ws.Range("N129").Value = -20106.0001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 1933
$ws.Range("I137").Value = 1908.5
$ws.Range("J137").Value = 1998.3334
$ws.Range("K137").Value = 5725.5
$ws.Range("L137").Value = 5995.0002
$ws.Range("M137").Value = -3175.5
$ws.Range("N137").Value = -11095.0002
$ws.Range("H138").Value = 1792.0605
$ws.Range("I138").Value = 1636.1852
$ws.Range("J138").Value = 2493.5
$ws.Range("K138").Value = 4908.5556
$ws.Range("L138").Value = 7480.5
$ws.Range("M138").Value = 231.4444000000003
$ws.Range("N138").Value = -17760.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4188.913
$ws.Range("I32").Value = 3833.9092
$ws.Range("K32").Value = 3833.9092
$ws.Range("M32").Value = -3546.9092
$ws.Range("H74").Value = 2112
$ws.Range("I74").Value = 2112
$ws.Range("K74").Value = 2112
$ws.Range("M74").Value = -1238
$ws.Range("H77").Value = 2112
$ws.Range("I77").Value = 2112
$ws.Range("K77").Value = 10560
$ws.Range("M77").Value = -6192
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 1218.2
$ws.Range("I132").Value = 1218.2
$ws.Range("K132").Value = 3654.6
$ws.Range("M132").Value = -1124.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 3833
$ws.Range("I24").Value = 4749.5
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 4749.5
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = -4514.5
$ws.Range("N24").Value = -2470
$ws.Range("H134").Value = 1899.875
$ws.Range("I134").Value = 1892.7142
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 5678.142599999999
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -3143.142599999999
$ws.Range("N134").Value = -10920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2120.111
$ws.Range("I31").Value = 1864.6364
$ws.Range("J31").Value = 2521.5715
$ws.Range("K31").Value = 1864.6364
$ws.Range("L31").Value = 2521.5715
$ws.Range("M31").Value = -1569.6364
$ws.Range("N31").Value = -3111.5715
$ws.Range("H34").Value = 2120.111
$ws.Range("I34").Value = 1864.6364
$ws.Range("J34").Value = 2521.5715
$ws.Range("K34").Value = 1864.6364
$ws.Range("L34").Value = 2521.5715
$ws.Range("M34").Value = -1662.6364
$ws.Range("N34").Value = -2925.5715
$ws.Range("H105").Value = 2883.48
$ws.Range("I105").Value = 2255.8572
$ws.Range("K105").Value = 2255.8572
$ws.Range("M105").Value = -508.8571999999999
$ws.Range("H122").Value = 2335.75
$ws.Range("I122").Value = 2443
$ws.Range("J122").Value = 1585
$ws.Range("K122").Value = 7329
$ws.Range("L122").Value = 4755
$ws.Range("M122").Value = -4879
$ws.Range("N122").Value = -9655
$ws.Range("H132").Value = 1840.96
$ws.Range("I132").Value = 1766.9546
$ws.Range("K132").Value = 5300.8638
$ws.Range("M132").Value = -2770.8638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 300
$ws.Range("M8").Value = -161
$ws.Range("H21").Value = 25
$ws.Range("I21").Value = 25
$ws.Range("K21").Value = 75
$ws.Range("M21").Value = 98
$ws.Range("H39").Value = 9447
$ws.Range("J39").Value = 9929.666999999999
$ws.Range("L39").Value = 29789.001
$ws.Range("N39").Value = -30377.001
$ws.Range("H56").Value = 20010
$ws.Range("I56").Value = 20010
$ws.Range("K56").Value = 20010
$ws.Range("M56").Value = -19480
$ws.Range("H116").Value = 906
$ws.Range("I116").Value = 999.5
$ws.Range("J116").Value = 719
$ws.Range("K116").Value = 2998.5
$ws.Range("L116").Value = 2157
$ws.Range("M116").Value = 443.5
$ws.Range("N116").Value = -9041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1666.6666
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -708
$ws.Range("H37").Value = 1666.6666
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -723
$ws.Range("H70").Value = 5933.75
$ws.Range("I70").Value = 5395.143
$ws.Range("J70").Value = 6687.8
$ws.Range("K70").Value = 5395.143
$ws.Range("L70").Value = 6687.8
$ws.Range("M70").Value = -5125.143
$ws.Range("N70").Value = -7227.8
$ws.Range("H73").Value = 5933.75
$ws.Range("I73").Value = 5395.143
$ws.Range("J73").Value = 6687.8
$ws.Range("K73").Value = 5395.143
$ws.Range("L73").Value = 6687.8
$ws.Range("M73").Value = -4459.143
$ws.Range("N73").Value = -8559.799999999999
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 2194
$ws.Range("I122").Value = 2194
$ws.Range("J122").Value = 2194
$ws.Range("K122").Value = 6582
$ws.Range("L122").Value = 6582
$ws.Range("M122").Value = -4132
$ws.Range("N122").Value = -11482
$ws.Range("H132").Value = 1965.8235
$ws.Range("I132").Value = 1855.6923
$ws.Range("J132").Value = 2323.75
$ws.Range("K132").Value = 5567.0769
$ws.Range("L132").Value = 6971.25
$ws.Range("M132").Value = -3037.0769
$ws.Range("N132").Value = -12031.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1712.5
$ws.Range("I16").Value = 1680.5555
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1680.5555
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1510.5555
$ws.Range("N16").Value = -2340
$ws.Range("H40").Value = 2629.7693
$ws.Range("I40").Value = 1899.5
$ws.Range("K40").Value = 1899.5
$ws.Range("M40").Value = -1763.5
$ws.Range("H122").Value = 7525.9443
$ws.Range("I122").Value = 12017
$ws.Range("J122").Value = 5280.4165
$ws.Range("K122").Value = 36051
$ws.Range("L122").Value = 15841.2495
$ws.Range("M122").Value = -33601
$ws.Range("N122").Value = -20741.2495
$ws.Range("H127").Value = 73748.5
$ws.Range("J127").Value = 73748.5
$ws.Range("L127").Value = 73748.5
$ws.Range("N127").Value = -83668.5
$ws.Range("H136").Value = 3832.6365
$ws.Range("I136").Value = 2906.111
$ws.Range("K136").Value = 8718.332999999999
$ws.Range("M136").Value = -6168.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1343.6666
$ws.Range("I113").Value = 1103
$ws.Range("J113").Value = 1825
$ws.Range("K113").Value = 3309
$ws.Range("L113").Value = 5475
$ws.Range("M113").Value = -1139
$ws.Range("N113").Value = -9815
$ws.Range("H119").Value = 135349
$ws.Range("J119").Value = 135349
$ws.Range("L119").Value = 135349
$ws.Range("N119").Value = -145025
$ws.Range("H132").Value = 2467.261
$ws.Range("I132").Value = 2561.7727
$ws.Range("K132").Value = 7685.3181
$ws.Range("M132").Value = -5155.3181
$ws.Range("H136").Value = 2457.9375
$ws.Range("I136").Value = 2371.8
$ws.Range("K136").Value = 7115.400000000001
$ws.Range("M136").Value = -4565.400000000001
